$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 772.625
$ws.Range("J17").Value = 762.73914
$ws.Range("L17").Value = 2288.21742
$ws.Range("N17").Value = -2624.21742
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1716
$ws.Range("H40").Value = 3233.2222
$ws.Range("J40").Value = 4299.8
$ws.Range("L40").Value = 4299.8
$ws.Range("N40").Value = -4649.8
$ws.Range("H69").Value = 185590.36
$ws.Range("I69").Value = 32499
$ws.Range("J69").Value = 200899.5
$ws.Range("K69").Value = 97497
$ws.Range("L69").Value = 602698.5
$ws.Range("M69").Value = -96623
$ws.Range("N69").Value = -604446.5
$ws.Range("H72").Value = 185590.36
$ws.Range("I72").Value = 32499
$ws.Range("J72").Value = 200899.5
$ws.Range("K72").Value = 292491
$ws.Range("L72").Value = 1808095.5
$ws.Range("M72").Value = -288123
$ws.Range("N72").Value = -1816831.5
$ws.Range("H80").Value = 7688.0586
$ws.Range("I80").Value = 1170.4286
$ws.Range("J80").Value = 12250.4
$ws.Range("K80").Value = 3511.2858
$ws.Range("L80").Value = 36751.2
$ws.Range("M80").Value = -2513.2858
$ws.Range("N80").Value = -38747.2
$ws.Range("H83").Value = 7688.0586
$ws.Range("I83").Value = 1170.4286
$ws.Range("J83").Value = 12250.4
$ws.Range("K83").Value = 10533.8574
$ws.Range("L83").Value = 110253.6
$ws.Range("M83").Value = -5541.857399999999
$ws.Range("N83").Value = -120237.6
$ws.Range("H86").Value = 2051.65
$ws.Range("I86").Value = 1771.4615
$ws.Range("J86").Value = 2572
$ws.Range("K86").Value = 1771.4615
$ws.Range("L86").Value = 2572
$ws.Range("M86").Value = -648.4614999999999
$ws.Range("N86").Value = -4818
$ws.Range("H89").Value = 2051.65
$ws.Range("I89").Value = 1771.4615
$ws.Range("J89").Value = 2572
$ws.Range("K89").Value = 8857.307499999999
$ws.Range("L89").Value = 12860
$ws.Range("M89").Value = -3241.307499999999
$ws.Range("N89").Value = -24092
$ws.Range("H98").Value = 1880.12
$ws.Range("I98").Value = 1130.3478
$ws.Range("K98").Value = 1130.3478
$ws.Range("M98").Value = 367.6522
$ws.Range("H101").Value = 591.1111
$ws.Range("I101").Value = 188.57143
$ws.Range("K101").Value = 565.71429
$ws.Range("M101").Value = 1056.28571
$ws.Range("H106").Value = 2460
$ws.Range("I106").Value = 2460
$ws.Range("K106").Value = 2460
$ws.Range("M106").Value = -1829
$ws.Range("H107").Value = 963.4666999999999
$ws.Range("I107").Value = 673.4
$ws.Range("J107").Value = 2413.8
$ws.Range("K107").Value = 673.4
$ws.Range("L107").Value = 2413.8
$ws.Range("M107").Value = 1246.6
$ws.Range("N107").Value = -6253.8
$ws.Range("H111").Value = 3308.25
$ws.Range("I111").Value = 3308.25
$ws.Range("K111").Value = 9924.75
$ws.Range("M111").Value = -6857.75
$ws.Range("H112").Value = 8915.4
$ws.Range("J112").Value = 9641.056
$ws.Range("L112").Value = 28923.168
$ws.Range("N112").Value = -31139.168
$ws.Range("H113").Value = 6381.593
$ws.Range("I113").Value = 6460.3335
$ws.Range("J113").Value = 6318.6
$ws.Range("K113").Value = 6460.3335
$ws.Range("L113").Value = 6318.6
$ws.Range("M113").Value = -3206.3335
$ws.Range("N113").Value = -12826.6
$ws.Range("H116").Value = 17694.592
$ws.Range("I116").Value = 25684.924
$ws.Range("J116").Value = 6153
$ws.Range("K116").Value = 25684.924
$ws.Range("L116").Value = 6153
$ws.Range("M116").Value = -22242.924
$ws.Range("N116").Value = -13037
$ws.Range("H122").Value = 1880.12
$ws.Range("I122").Value = 1130.3478
$ws.Range("K122").Value = 3391.0434
$ws.Range("M122").Value = -941.0434
$ws.Range("H125").Value = 47503.715
$ws.Range("I125").Value = 104676.664
$ws.Range("J125").Value = 4624
$ws.Range("K125").Value = 942089.976
$ws.Range("L125").Value = 41616
$ws.Range("M125").Value = -939629.976
$ws.Range("N125").Value = -46536
$ws.Range("H132").Value = 22430.146
$ws.Range("I132").Value = 24917.5
$ws.Range("J132").Value = 3775
$ws.Range("K132").Value = 74752.5
$ws.Range("L132").Value = 11325
$ws.Range("M132").Value = -72222.5
$ws.Range("N132").Value = -16385
$ws.Range("H137").Value = 24941.52
$ws.Range("I137").Value = 18554.117
$ws.Range("J137").Value = 38514.75
$ws.Range("K137").Value = 55662.351
$ws.Range("L137").Value = 115544.25
$ws.Range("M137").Value = -53112.351
$ws.Range("N137").Value = -120644.25
$ws.Range("H138").Value = 25343.303
$ws.Range("I138").Value = 1596.5454
$ws.Range("J138").Value = 103707.6
$ws.Range("K138").Value = 4789.6362
$ws.Range("L138").Value = 311122.8
$ws.Range("M138").Value = 350.3638000000001
$ws.Range("N138").Value = -321402.8
$ws.Range("H141").Value = 2177.1
$ws.Range("I141").Value = 2177.1
$ws.Range("K141").Value = 6531.299999999999
$ws.Range("M141").Value = -1351.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2731.5334
$ws.Range("I2").Value = 2834.04
$ws.Range("J2").Value = 2219
$ws.Range("K2").Value = 2834.04
$ws.Range("L2").Value = 2219
$ws.Range("M2").Value = -2721.04
$ws.Range("N2").Value = -2445
$ws.Range("H32").Value = 20395.426
$ws.Range("I32").Value = 22783.521
$ws.Range("K32").Value = 22783.521
$ws.Range("M32").Value = -22496.521
$ws.Range("H43").Value = 20767.727
$ws.Range("J43").Value = 19372.857
$ws.Range("L43").Value = 19372.857
$ws.Range("N43").Value = -19998.857
$ws.Range("H45").Value = 3836.1177
$ws.Range("I45").Value = 2128.6
$ws.Range("J45").Value = 6275.4287
$ws.Range("K45").Value = 2128.6
$ws.Range("L45").Value = 6275.4287
$ws.Range("M45").Value = -1751.6
$ws.Range("N45").Value = -7029.4287
$ws.Range("H61").Value = 5298.7915
$ws.Range("I61").Value = 1007.7727
$ws.Range("J61").Value = 52500
$ws.Range("K61").Value = 1007.7727
$ws.Range("L61").Value = 52500
$ws.Range("M61").Value = -795.7727
$ws.Range("N61").Value = -52924
$ws.Range("H74").Value = 360784.6
$ws.Range("I74").Value = 667550.7
$ws.Range("J74").Value = 15672.75
$ws.Range("K74").Value = 667550.7
$ws.Range("L74").Value = 15672.75
$ws.Range("M74").Value = -666676.7
$ws.Range("N74").Value = -17420.75
$ws.Range("H77").Value = 360784.6
$ws.Range("I77").Value = 667550.7
$ws.Range("J77").Value = 15672.75
$ws.Range("K77").Value = 3337753.5
$ws.Range("L77").Value = 78363.75
$ws.Range("M77").Value = -3333385.5
$ws.Range("N77").Value = -87099.75
$ws.Range("H97").Value = 1432.2667
$ws.Range("I97").Value = 1123
$ws.Range("J97").Value = 2978.6
$ws.Range("K97").Value = 1123
$ws.Range("L97").Value = 2978.6
$ws.Range("M97").Value = -627
$ws.Range("N97").Value = -3970.6
$ws.Range("H102").Value = 3530
$ws.Range("I102").Value = 3353.8
$ws.Range("J102").Value = 4411
$ws.Range("K102").Value = 3353.8
$ws.Range("L102").Value = 4411
$ws.Range("M102").Value = -1731.8
$ws.Range("N102").Value = -7655
$ws.Range("H109").Value = 99877
$ws.Range("J109").Value = 99877
$ws.Range("L109").Value = 99877
$ws.Range("N109").Value = -102651
$ws.Range("H116").Value = 2731.5334
$ws.Range("I116").Value = 2834.04
$ws.Range("J116").Value = 2219
$ws.Range("K116").Value = 2834.04
$ws.Range("L116").Value = 2219
$ws.Range("M116").Value = -540.04
$ws.Range("N116").Value = -6807
$ws.Range("H122").Value = 1796.32
$ws.Range("I122").Value = 1540.6818
$ws.Range("K122").Value = 4622.0454
$ws.Range("M122").Value = -2172.0454
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""
$ws.Range("H132").Value = 2124.5
$ws.Range("I132").Value = 1199.2
$ws.Range("K132").Value = 3597.6
$ws.Range("M132").Value = -1067.6
$ws.Range("H136").Value = 5298.7915
$ws.Range("I136").Value = 1007.7727
$ws.Range("J136").Value = 52500
$ws.Range("K136").Value = 3023.3181
$ws.Range("L136").Value = 157500
$ws.Range("M136").Value = -473.3181
$ws.Range("N136").Value = -162600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2731.5334
$ws.Range("I3").Value = 2834.04
$ws.Range("J3").Value = 2219
$ws.Range("K3").Value = 2834.04
$ws.Range("L3").Value = 2219
$ws.Range("M3").Value = -2720.04
$ws.Range("N3").Value = -2447
$ws.Range("H22").Value = 1000.3333
$ws.Range("I22").Value = 1000.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1000.5
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -827.5
$ws.Range("N22").Value = -1346
$ws.Range("H40").Value = 32250
$ws.Range("J40").Value = 32250
$ws.Range("L40").Value = 32250
$ws.Range("N40").Value = -32780
$ws.Range("H80").Value = 764.17645
$ws.Range("I80").Value = 855.55554
$ws.Range("K80").Value = 855.55554
$ws.Range("M80").Value = 142.44446
$ws.Range("H83").Value = 764.17645
$ws.Range("I83").Value = 855.55554
$ws.Range("K83").Value = 4277.7777
$ws.Range("M83").Value = 714.2223000000004
$ws.Range("H86").Value = 3166.3333
$ws.Range("I86").Value = 2999
$ws.Range("K86").Value = 2999
$ws.Range("M86").Value = -1876
$ws.Range("H89").Value = 3166.3333
$ws.Range("I89").Value = 2999
$ws.Range("K89").Value = 14995
$ws.Range("M89").Value = -9379
$ws.Range("H94").Value = 8114.643
$ws.Range("I94").Value = 8967.166999999999
$ws.Range("K94").Value = 8967.166999999999
$ws.Range("M94").Value = -8516.166999999999
$ws.Range("H96").Value = 9425
$ws.Range("I96").Value = 9425
$ws.Range("K96").Value = 9425
$ws.Range("M96").Value = -6679
$ws.Range("H99").Value = 1025.1428
$ws.Range("I99").Value = 1047.1666
$ws.Range("J99").Value = 893
$ws.Range("K99").Value = 1047.1666
$ws.Range("L99").Value = 893
$ws.Range("M99").Value = 450.8334
$ws.Range("N99").Value = -3889
$ws.Range("H107").Value = 5122.8184
$ws.Range("J107").Value = 5288.6665
$ws.Range("L107").Value = 5288.6665
$ws.Range("N107").Value = -9128.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1237.5
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650
$ws.Range("H31").Value = 6251231
$ws.Range("I31").Value = 11112022
$ws.Range("J31").Value = 1642.8572
$ws.Range("K31").Value = 11112022
$ws.Range("L31").Value = 1642.8572
$ws.Range("M31").Value = -11111727
$ws.Range("N31").Value = -2232.8572
$ws.Range("H34").Value = 6251231
$ws.Range("I34").Value = 11112022
$ws.Range("J34").Value = 1642.8572
$ws.Range("K34").Value = 11112022
$ws.Range("L34").Value = 1642.8572
$ws.Range("M34").Value = -11111820
$ws.Range("N34").Value = -2046.8572
$ws.Range("H55").Value = 47000
$ws.Range("J55").Value = 47000
$ws.Range("L55").Value = 47000
$ws.Range("N55").Value = -47630
$ws.Range("H58").Value = 1277.9524
$ws.Range("I58").Value = 1048.1875
$ws.Range("J58").Value = 2013.2
$ws.Range("K58").Value = 1048.1875
$ws.Range("L58").Value = 2013.2
$ws.Range("M58").Value = -845.1875
$ws.Range("N58").Value = -2419.2
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -81240
$ws.Range("H99").Value = 7432.3335
$ws.Range("I99").Value = 6913.8335
$ws.Range("J99").Value = 8469.333000000001
$ws.Range("K99").Value = 6913.8335
$ws.Range("L99").Value = 8469.333000000001
$ws.Range("M99").Value = -5415.8335
$ws.Range("N99").Value = -11465.333
$ws.Range("H105").Value = 50217.125
$ws.Range("I105").Value = 77947.39999999999
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 77947.39999999999
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -76200.39999999999
$ws.Range("N105").Value = -7494
$ws.Range("H107").Value = 850.8095
$ws.Range("I107").Value = 703.13336
$ws.Range("J107").Value = 1220
$ws.Range("K107").Value = 703.13336
$ws.Range("L107").Value = 1220
$ws.Range("M107").Value = 1216.86664
$ws.Range("N107").Value = -5060
$ws.Range("H122").Value = 2121.8
$ws.Range("I122").Value = 2121.8
$ws.Range("K122").Value = 6365.400000000001
$ws.Range("M122").Value = -3915.400000000001
$ws.Range("H126").Value = 7432.3335
$ws.Range("I126").Value = 6913.8335
$ws.Range("J126").Value = 8469.333000000001
$ws.Range("K126").Value = 20741.5005
$ws.Range("L126").Value = 25407.999
$ws.Range("M126").Value = -18271.5005
$ws.Range("N126").Value = -30347.999
$ws.Range("H132").Value = 126711.25
$ws.Range("I132").Value = 167282.5
$ws.Range("J132").Value = 4997.5
$ws.Range("K132").Value = 501847.5
$ws.Range("L132").Value = 14992.5
$ws.Range("M132").Value = -499317.5
$ws.Range("N132").Value = -20052.5
$ws.Range("H134").Value = 2620.3572
$ws.Range("I134").Value = 2419.625
$ws.Range("J134").Value = 3824.75
$ws.Range("K134").Value = 7258.875
$ws.Range("L134").Value = 11474.25
$ws.Range("M134").Value = -4723.875
$ws.Range("N134").Value = -16544.25
$ws.Range("H136").Value = 1277.9524
$ws.Range("I136").Value = 1048.1875
$ws.Range("J136").Value = 2013.2
$ws.Range("K136").Value = 3144.5625
$ws.Range("L136").Value = 6039.6
$ws.Range("M136").Value = -594.5625
$ws.Range("N136").Value = -11139.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79637656
$ws.Range("I4").Value = 82331100
$ws.Range("K4").Value = 246993300
$ws.Range("M4").Value = -246993188
$ws.Range("H5").Value = 784.6667
$ws.Range("I5").Value = 734.75
$ws.Range("K5").Value = 2204.25
$ws.Range("M5").Value = -2092.25
$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665
$ws.Range("H69").Value = 4863.1577
$ws.Range("I69").Value = 2900
$ws.Range("J69").Value = 4972.222
$ws.Range("K69").Value = 8700
$ws.Range("L69").Value = 14916.666
$ws.Range("M69").Value = -7889
$ws.Range("N69").Value = -16538.666
$ws.Range("H72").Value = 4863.1577
$ws.Range("I72").Value = 2900
$ws.Range("J72").Value = 4972.222
$ws.Range("K72").Value = 26100
$ws.Range("L72").Value = 44749.998
$ws.Range("M72").Value = -22044
$ws.Range("N72").Value = -52861.998
$ws.Range("H98").Value = 1901
$ws.Range("I98").Value = 1351.5
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 4054.5
$ws.Range("L98").Value = 9000
$ws.Range("M98").Value = -2556.5
$ws.Range("N98").Value = -11996
$ws.Range("H107").Value = 2006.2307
$ws.Range("I107").Value = 2871.7144
$ws.Range("J107").Value = 996.5
$ws.Range("K107").Value = 8615.143199999999
$ws.Range("L107").Value = 2989.5
$ws.Range("M107").Value = -6695.143199999999
$ws.Range("N107").Value = -6829.5
$ws.Range("H110").Value = 11512.5
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H114").Value = 760.1429000000001
$ws.Range("J114").Value = 818.6
$ws.Range("L114").Value = 2455.8
$ws.Range("N114").Value = -8963.799999999999
$ws.Range("H131").Value = 111245.44
$ws.Range("I131").Value = 390163.72
$ws.Range("K131").Value = 1170491.16
$ws.Range("M131").Value = -1165451.16
$ws.Range("H135").Value = 784.6667
$ws.Range("I135").Value = 734.75
$ws.Range("K135").Value = 6612.75
$ws.Range("M135").Value = -4077.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22555.555
$ws.Range("J15").Value = 22555.555
$ws.Range("L15").Value = 22555.555
$ws.Range("N15").Value = -23131.555
$ws.Range("H81").Value = 22555.555
$ws.Range("J81").Value = 22555.555
$ws.Range("L81").Value = 22555.555
$ws.Range("N81").Value = -24551.555
$ws.Range("H84").Value = 22555.555
$ws.Range("J84").Value = 22555.555
$ws.Range("L84").Value = 67666.66500000001
$ws.Range("N84").Value = -77650.66500000001
$ws.Range("H97").Value = 1191.6
$ws.Range("J97").Value = 1473.4286
$ws.Range("L97").Value = 1473.4286
$ws.Range("N97").Value = -2465.4286
$ws.Range("H113").Value = 3147.75
$ws.Range("I113").Value = 3012.6924
$ws.Range("K113").Value = 3012.6924
$ws.Range("M113").Value = -842.6923999999999
$ws.Range("H114").Value = 74979
$ws.Range("J114").Value = 74979
$ws.Range("L114").Value = 74979
$ws.Range("N114").Value = -83657
$ws.Range("H122").Value = 3907.2144
$ws.Range("I122").Value = 3708.5557
$ws.Range("J122").Value = 4264.8
$ws.Range("K122").Value = 11125.6671
$ws.Range("L122").Value = 12794.4
$ws.Range("M122").Value = -8675.667099999999
$ws.Range("N122").Value = -17694.4
$ws.Range("H132").Value = 2327.2856
$ws.Range("I132").Value = 2092.3333
$ws.Range("J132").Value = 2503.5
$ws.Range("K132").Value = 6276.999899999999
$ws.Range("L132").Value = 7510.5
$ws.Range("M132").Value = -3746.999899999999
$ws.Range("N132").Value = -12570.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3624.8696
$ws.Range("I7").Value = 3624.8696
$ws.Range("K7").Value = 3624.8696
$ws.Range("M7").Value = -3512.8696
$ws.Range("H16").Value = 2034.35
$ws.Range("I16").Value = 1299.375
$ws.Range("J16").Value = 4974.25
$ws.Range("K16").Value = 1299.375
$ws.Range("L16").Value = 4974.25
$ws.Range("M16").Value = -1129.375
$ws.Range("N16").Value = -5314.25
$ws.Range("H22").Value = 1558.8823
$ws.Range("I22").Value = 1230.1
$ws.Range("J22").Value = 2028.5714
$ws.Range("K22").Value = 1230.1
$ws.Range("L22").Value = 2028.5714
$ws.Range("M22").Value = -935.0999999999999
$ws.Range("N22").Value = -2618.5714
$ws.Range("H27").Value = 1558.8823
$ws.Range("I27").Value = 1230.1
$ws.Range("J27").Value = 2028.5714
$ws.Range("K27").Value = 1230.1
$ws.Range("L27").Value = 2028.5714
$ws.Range("M27").Value = -1123.1
$ws.Range("N27").Value = -2242.5714
$ws.Range("H40").Value = 3175.3333
$ws.Range("I40").Value = 3048.5386
$ws.Range("J40").Value = 3999.5
$ws.Range("K40").Value = 3048.5386
$ws.Range("L40").Value = 3999.5
$ws.Range("M40").Value = -2912.5386
$ws.Range("N40").Value = -4271.5
$ws.Range("H46").Value = 1979.2285
$ws.Range("I46").Value = 1206.375
$ws.Range("J46").Value = 2630.0527
$ws.Range("K46").Value = 1206.375
$ws.Range("L46").Value = 2630.0527
$ws.Range("M46").Value = -1018.375
$ws.Range("N46").Value = -3006.0527
$ws.Range("H61").Value = 1933.8
$ws.Range("I61").Value = 1667.25
$ws.Range("K61").Value = 1667.25
$ws.Range("M61").Value = -1465.25
$ws.Range("H68").Value = 3937
$ws.Range("J68").Value = 5166.1665
$ws.Range("L68").Value = 5166.1665
$ws.Range("N68").Value = -6664.1665
$ws.Range("H71").Value = 3937
$ws.Range("J71").Value = 5166.1665
$ws.Range("L71").Value = 25830.8325
$ws.Range("N71").Value = -33318.8325
$ws.Range("H74").Value = 149999.8
$ws.Range("I74").Value = 149999
$ws.Range("K74").Value = 149999
$ws.Range("M74").Value = -149001
$ws.Range("H77").Value = 149999.8
$ws.Range("I77").Value = 149999
$ws.Range("K77").Value = 449997
$ws.Range("M77").Value = -445005
$ws.Range("H87").Value = 98571.42999999999
$ws.Range("J87").Value = 98571.42999999999
$ws.Range("L87").Value = 98571.42999999999
$ws.Range("N87").Value = -100817.43
$ws.Range("H88").Value = 94285.57000000001
$ws.Range("I88").Value = 55000
$ws.Range("J88").Value = 100833.164
$ws.Range("K88").Value = 55000
$ws.Range("L88").Value = 100833.164
$ws.Range("M88").Value = -54572
$ws.Range("N88").Value = -101689.164
$ws.Range("H90").Value = 98571.42999999999
$ws.Range("J90").Value = 98571.42999999999
$ws.Range("L90").Value = 295714.29
$ws.Range("N90").Value = -306946.29
$ws.Range("H91").Value = 94285.57000000001
$ws.Range("I91").Value = 55000
$ws.Range("J91").Value = 100833.164
$ws.Range("K91").Value = 55000
$ws.Range("L91").Value = 100833.164
$ws.Range("M91").Value = -53518
$ws.Range("N91").Value = -103797.164
$ws.Range("H93").Value = 1549
$ws.Range("I93").Value = 1319
$ws.Range("K93").Value = 1319
$ws.Range("M93").Value = -71
$ws.Range("H100").Value = 2753.8235
$ws.Range("I100").Value = 2471.3635
$ws.Range("K100").Value = 2471.3635
$ws.Range("M100").Value = -1930.3635
$ws.Range("H113").Value = 1933.8
$ws.Range("I113").Value = 1667.25
$ws.Range("K113").Value = 1667.25
$ws.Range("M113").Value = 502.75
$ws.Range("H122").Value = 3142.4666
$ws.Range("I122").Value = 3489.8
$ws.Range("K122").Value = 10469.4
$ws.Range("M122").Value = -8019.400000000001
$ws.Range("H126").Value = 3624.8696
$ws.Range("I126").Value = 3624.8696
$ws.Range("K126").Value = 10874.6088
$ws.Range("M126").Value = -8404.6088
$ws.Range("H132").Value = 3116.027
$ws.Range("I132").Value = 2827.3794
$ws.Range("J132").Value = 4162.375
$ws.Range("K132").Value = 8482.138199999999
$ws.Range("L132").Value = 12487.125
$ws.Range("M132").Value = -5952.138199999999
$ws.Range("N132").Value = -17547.125
$ws.Range("H136").Value = 4755.8184
$ws.Range("I136").Value = 4735.6665
$ws.Range("J136").Value = 4763.375
$ws.Range("K136").Value = 14206.9995
$ws.Range("L136").Value = 14290.125
$ws.Range("M136").Value = -11656.9995
$ws.Range("N136").Value = -19390.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18416.375
$ws.Range("J45").Value = 18701.715
$ws.Range("L45").Value = 18701.715
$ws.Range("N45").Value = -19683.715
$ws.Range("H75").Value = 36133.332
$ws.Range("J75").Value = 36133.332
$ws.Range("L75").Value = 36133.332
$ws.Range("N75").Value = -38005.332
$ws.Range("H78").Value = 36133.332
$ws.Range("J78").Value = 36133.332
$ws.Range("L78").Value = 108399.996
$ws.Range("N78").Value = -117759.996
$ws.Range("H81").Value = 5767.3687
$ws.Range("I81").Value = 6608.6665
$ws.Range("J81").Value = 2612.5
$ws.Range("K81").Value = 13217.333
$ws.Range("L81").Value = 5225
$ws.Range("M81").Value = -12156.333
$ws.Range("N81").Value = -7347
$ws.Range("H84").Value = 5767.3687
$ws.Range("I84").Value = 6608.6665
$ws.Range("J84").Value = 2612.5
$ws.Range("K84").Value = 66086.66500000001
$ws.Range("L84").Value = 26125
$ws.Range("M84").Value = -60782.66500000001
$ws.Range("N84").Value = -36733
$ws.Range("H86").Value = 67250
$ws.Range("J86").Value = 67250
$ws.Range("L86").Value = 67250
$ws.Range("N86").Value = -69496
$ws.Range("H89").Value = 67250
$ws.Range("J89").Value = 67250
$ws.Range("L89").Value = 336250
$ws.Range("N89").Value = -347482
$ws.Range("H104").Value = 29329.334
$ws.Range("J104").Value = 29329.334
$ws.Range("L104").Value = 29329.334
$ws.Range("N104").Value = -36317.334
$ws.Range("H107").Value = 2133
$ws.Range("I107").Value = 2133
$ws.Range("K107").Value = 6399
$ws.Range("M107").Value = -4479
$ws.Range("H113").Value = 1553.45
$ws.Range("I113").Value = 1322.0834
$ws.Range("J113").Value = 1900.5
$ws.Range("K113").Value = 3966.2502
$ws.Range("L113").Value = 5701.5
$ws.Range("M113").Value = -1796.2502
$ws.Range("N113").Value = -10041.5
$ws.Range("H122").Value = 10901184
$ws.Range("I122").Value = 12857283
$ws.Range("J122").Value = 2917.2856
$ws.Range("K122").Value = 38571849
$ws.Range("L122").Value = 8751.856800000001
$ws.Range("M122").Value = -38569399
$ws.Range("N122").Value = -13651.8568
$ws.Range("H124").Value = 59500
$ws.Range("J124").Value = 59500
$ws.Range("L124").Value = 59500
$ws.Range("N124").Value = -69320
$ws.Range("H126").Value = 359965.28
$ws.Range("I126").Value = 2517.2222
$ws.Range("J126").Value = 1003371.8
$ws.Range("K126").Value = 7551.6666
$ws.Range("L126").Value = 3010115.4
$ws.Range("M126").Value = -5081.6666
$ws.Range("N126").Value = -3015055.4
$ws.Range("H132").Value = 7598655
$ws.Range("I132").Value = 8645958
$ws.Range("K132").Value = 25937874
$ws.Range("M132").Value = -25935344
$ws.Range("H136").Value = 20185.943
$ws.Range("I136").Value = 24232.5
$ws.Range("J136").Value = 3999.7144
$ws.Range("K136").Value = 72697.5
$ws.Range("L136").Value = 11999.1432
$ws.Range("M136").Value = -70147.5
$ws.Range("N136").Value = -17099.1432
